$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-77 down to 42-78
$ws.Rows.Item(41).EntireRow.Insert()

# Text columns must stay as text (matches the rest of the sheet), so force
# text format before assigning values that would otherwise be auto-typed
# (numbers / dates) by Excel.
$textRange = $ws.Range("A41:H41")
$textRange.NumberFormat = "@"
$textRange2 = $ws.Range("J41:L41")
$textRange2.NumberFormat = "@"
$textRange3 = $ws.Range("O41:P41")
$textRange3.NumberFormat = "@"

# Populate the newly inserted row 41 with the new record
$ws.Range("A41").Value = "-603"
$ws.Range("B41").Value = "9/22/2025"
$ws.Range("C41").Value = "ANCHORENA, TOMAS MANUEL DE, DR. 821"
$ws.Range("D41").Value = "3"
$ws.Range("E41").Value = "809910086"
$ws.Range("F41").Value = "PEBCOM"
$ws.Range("G41").Value = "Pendiente"
$ws.Range("H41").Value = "Columna chocada pendiente para instalar un corporativo"
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = "Cambio"
$ws.Range("K41").Value = "Sin equipos"
$ws.Range("L41").Value = "Pasante"
$ws.Range("M41").Value = -58.408551
$ws.Range("N41").Value = -34.599265
$ws.Range("O41").Value = "Almagro"
$ws.Range("P41").Value = "Capital Sur"
